$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032640661117163
$ws.Range("D2").Value = 1.04389119359344
$ws.Range("E2").Value = 1.032061870338741
$ws.Range("F2").Value = 1.054845114770408
$ws.Range("I2").Value = 1.041792787156781
$ws.Range("J2").Value = 1.037769767682524
$ws.Range("K2").Value = 1.046663782943812
$ws.Range("L2").Value = 1.034868186670981
$ws.Range("M2").Value = 1.057587204250136
$ws.Range("N2").Value = 1.016571308471184

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03343642842062
$ws.Range("D3").Value = 1.044543721015948
$ws.Range("E3").Value = 1.032734065265371
$ws.Range("F3").Value = 1.055666363058192
$ws.Range("I3").Value = 1.042012014437832
$ws.Range("J3").Value = 1.038208695475633
$ws.Range("K3").Value = 1.047127851775967
$ws.Range("L3").Value = 1.035349440994537
$ws.Range("M3").Value = 1.058221756177087
$ws.Range("N3").Value = 1.016717575964628

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033951959870905
$ws.Range("D4").Value = 1.04496644895903
$ws.Range("E4").Value = 1.033169937314351
$ws.Range("F4").Value = 1.056198643299491
$ws.Range("I4").Value = 1.042152871726409
$ws.Range("J4").Value = 1.038492649858487
$ws.Range("K4").Value = 1.047427947920108
$ws.Range("L4").Value = 1.035661080081816
$ws.Range("M4").Value = 1.058632578439082
$ws.Range("N4").Value = 1.016812173786673

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034168835183674
$ws.Range("D5").Value = 1.045144281383541
$ws.Range("E5").Value = 1.033353395564738
$ws.Range("F5").Value = 1.056422621967682
$ws.Range("I5").Value = 1.042211848468954
$ws.Range("J5").Value = 1.038612008333313
$ws.Range("K5").Value = 1.047554062055988
$ws.Range("L5").Value = 1.035792148094232
$ws.Range("M5").Value = 1.0588053401764
$ws.Range("N5").Value = 1.0168519309882

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03420525800667
$ws.Range("D6").Value = 1.045174147071354
$ws.Range("E6").Value = 1.033384211719896
$ws.Range("F6").Value = 1.056460241112017
$ws.Range("I6").Value = 1.0422217368442
$ws.Range("J6").Value = 1.038632048173902
$ws.Range("K6").Value = 1.047575234419145
$ws.Range("L6").Value = 1.035814158171071
$ws.Range("M6").Value = 1.058834350634314
$ws.Range("N6").Value = 1.016858605696327

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033954857198412
$ws.Range("D7").Value = 1.044968824704878
$ws.Range("E7").Value = 1.033172387839666
$ws.Range("F7").Value = 1.056201635298496
$ws.Range("I7").Value = 1.042153660719329
$ws.Range("J7").Value = 1.038494244796544
$ws.Range("K7").Value = 1.047429633246999
$ws.Range("L7").Value = 1.035662831204902
$ws.Range("M7").Value = 1.058634886687506
$ws.Range("N7").Value = 1.016812705070993

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03290946562692
$ws.Range("D8").Value = 1.044111613825814
$ws.Range("E8").Value = 1.032288850760503
$ws.Range("F8").Value = 1.055122476503039
$ws.Range("I8").Value = 1.041867082038048
$ws.Range("J8").Value = 1.037918117295275
$ws.Range("K8").Value = 1.046820655135331
$ws.Range("L8").Value = 1.035030779245561
$ws.Range("M8").Value = 1.057801606464117
$ws.Range("N8").Value = 1.016620749713474

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03107215437854
$ws.Range("D9").Value = 1.0426050041001
$ws.Range("E9").Value = 1.030739051908677
$ws.Range("F9").Value = 1.053227674134899
$ws.Range("I9").Value = 1.041354493541941
$ws.Range("J9").Value = 1.036902494719468
$ws.Range("K9").Value = 1.045746184985591
$ws.Range("L9").Value = 1.033918887628848
$ws.Range("M9").Value = 1.056335058532505
$ws.Range("N9").Value = 1.016282160207294

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029850614632654
$ws.Range("D10").Value = 1.041603335607996
$ws.Range("E10").Value = 1.02971074327645
$ws.Range("F10").Value = 1.051969176701077
$ws.Range("I10").Value = 1.041007710263022
$ws.Range("J10").Value = 1.036225212429545
$ws.Range("K10").Value = 1.045029032862594
$ws.Range("L10").Value = 1.033178964154312
$ws.Range("M10").Value = 1.055358673038432
$ws.Range("N10").Value = 1.016056233225194

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029322487596133
$ws.Range("D11").Value = 1.041170275253848
$ws.Range("E11").Value = 1.02926665814097
$ws.Range("F11").Value = 1.051425375238129
$ws.Range("I11").Value = 1.040856360713405
$ws.Range("J11").Value = 1.03593191033463
$ws.Range("K11").Value = 1.044718316972565
$ws.Range("L11").Value = 1.032858904125108
$ws.Range("M11").Value = 1.054936220783909
$ws.Range("N11").Value = 1.015958362612499

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029126440848492
$ws.Range("D12").Value = 1.041009519749126
$ws.Range("E12").Value = 1.029101884268328
$ws.Range("F12").Value = 1.051223555981023
$ws.Range("I12").Value = 1.040799964827492
$ws.Range("J12").Value = 1.035822961014343
$ws.Range("K12").Value = 1.044602877004252
$ws.Range("L12").Value = 1.032740070853193
$ws.Range("M12").Value = 1.054779354427506
$ws.Range("N12").Value = 1.015922003160424

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029168487939775
$ws.Range("D13").Value = 1.041043997686197
$ws.Range("E13").Value = 1.02913722068662
$ws.Range("F13").Value = 1.051266839042423
$ws.Range("I13").Value = 1.040812069973851
$ws.Range("J13").Value = 1.035846331167254
$ws.Range("K13").Value = 1.044627640434256
$ws.Range("L13").Value = 1.03276555865012
$ws.Range("M13").Value = 1.054813000447103
$ws.Range("N13").Value = 1.015929802647942

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029306279788231
$ws.Range("D14").Value = 1.041156985053768
$ws.Range("E14").Value = 1.029253034216037
$ws.Range("F14").Value = 1.051408689257707
$ws.Range("I14").Value = 1.040851702636421
$ws.Range("J14").Value = 1.035922904621091
$ws.Range("K14").Value = 1.044708775196012
$ws.Range("L14").Value = 1.032849080275513
$ws.Range("M14").Value = 1.054923253106802
$ws.Range("N14").Value = 1.015955357247475

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029391194265895
$ws.Range("D15").Value = 1.041226613925594
$ws.Range("E15").Value = 1.029324414547033
$ws.Range("F15").Value = 1.051496110820375
$ws.Range("I15").Value = 1.040876098077366
$ws.Range("J15").Value = 1.035970083580368
$ws.Range("K15").Value = 1.044758761564767
$ws.Range("L15").Value = 1.032900547542676
$ws.Range("M15").Value = 1.054991190249576
$ws.Range("N15").Value = 1.01597110150549

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029885681855986
$ws.Range("D16").Value = 1.041632090650951
$ws.Range("E16").Value = 1.029740240756668
$ws.Range("F16").Value = 1.052005291132684
$ws.Range("I16").Value = 1.041017729834754
$ws.Range("J16").Value = 1.036244677310327
$ws.Range("K16").Value = 1.045049650292669
$ws.Range("L16").Value = 1.033200212595949
$ws.Range("M16").Value = 1.055386716888776
$ws.Range("N16").Value = 1.016062727715292

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030196078377453
$ws.Range("D17").Value = 1.041886615882305
$ws.Range("E17").Value = 1.030001394654893
$ws.Range("F17").Value = 1.052324992183516
$ws.Range("I17").Value = 1.041106253687559
$ws.Range("J17").Value = 1.036416914415173
$ws.Range("K17").Value = 1.045232068780581
$ws.Range("L17").Value = 1.033388274391886
$ws.Range("M17").Value = 1.05563490966519
$ws.Range("N17").Value = 1.016120191307121

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030377205231065
$ws.Range("D18").Value = 1.042035140535595
$ws.Range("E18").Value = 1.030153834959307
$ws.Range("F18").Value = 1.052511577885678
$ws.Range("I18").Value = 1.041157773266912
$ws.Range("J18").Value = 1.036517373998445
$ws.Range("K18").Value = 1.04533845252095
$ws.Range("L18").Value = 1.033497999565417
$ws.Range("M18").Value = 1.055779707947626
$ws.Range("N18").Value = 1.016153704692742

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030438977951437
$ws.Range("D19").Value = 1.042085794467517
$ws.Range("E19").Value = 1.030205832382168
$ws.Range("F19").Value = 1.052575217295092
$ws.Range("I19").Value = 1.041175320598783
$ws.Range("J19").Value = 1.036551627502833
$ws.Range("K19").Value = 1.045374723551235
$ws.Range("L19").Value = 1.033535418382749
$ws.Range("M19").Value = 1.05582908573183
$ws.Range("N19").Value = 1.016165131171995

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030162767706631
$ws.Range("D20").Value = 1.041859301059514
$ws.Range("E20").Value = 1.029973363543359
$ws.Range("F20").Value = 1.052290679946269
$ws.Range("I20").Value = 1.041096767792852
$ws.Range("J20").Value = 1.036398435346744
$ws.Range("K20").Value = 1.045212498842006
$ws.Range("L20").Value = 1.033368093830046
$ws.Range("M20").Value = 1.055608277644182
$ws.Range("N20").Value = 1.016114026433941

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029265700109299
$ws.Range("D21").Value = 1.041123710248155
$ws.Range("E21").Value = 1.029218925058021
$ws.Range("F21").Value = 1.051366913117213
$ws.Range("I21").Value = 1.040840036709662
$ws.Range("J21").Value = 1.035900355743643
$ws.Range("K21").Value = 1.044684883744071
$ws.Range("L21").Value = 1.032824483811625
$ws.Range("M21").Value = 1.054890785018297
$ws.Range("N21").Value = 1.015947832215272

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028702391496351
$ws.Range("D22").Value = 1.040661808474368
$ws.Range("E22").Value = 1.028745616929496
$ws.Range("F22").Value = 1.050787105170305
$ws.Range("I22").Value = 1.040677590706756
$ws.Range("J22").Value = 1.035587171623793
$ws.Range("K22").Value = 1.044353000014185
$ws.Range("L22").Value = 1.032482991980481
$ws.Range("M22").Value = 1.054439966152506
$ws.Range("N22").Value = 1.015843305128978

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029000943828981
$ws.Range("D23").Value = 1.040906614457493
$ws.Range("E23").Value = 1.028996427552011
$ws.Range("F23").Value = 1.051094376752017
$ws.Range("I23").Value = 1.040763803649484
$ws.Range("J23").Value = 1.035753198073325
$ws.Range("K23").Value = 1.044528951709631
$ws.Range("L23").Value = 1.032663994595491
$ws.Range("M23").Value = 1.054678924996671
$ws.Range("N23").Value = 1.015898719997947

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030177819124727
$ws.Range("D24").Value = 1.041871643253472
$ws.Range("E24").Value = 1.029986029245163
$ws.Range("F24").Value = 1.052306183831301
$ws.Range("I24").Value = 1.041101054415442
$ws.Range("J24").Value = 1.036406785253035
$ws.Range("K24").Value = 1.045221341710291
$ws.Range("L24").Value = 1.03337721245887
$ws.Range("M24").Value = 1.055620311410594
$ws.Range("N24").Value = 1.016116812087667

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031546563404294
$ws.Range("D25").Value = 1.042994024446681
$ws.Range("E25").Value = 1.031138858297075
$ws.Range("F25").Value = 1.053716705849159
$ws.Range("I25").Value = 1.041487905232785
$ws.Range("J25").Value = 1.037165098461402
$ws.Range("K25").Value = 1.046024115274688
$ws.Range("L25").Value = 1.034206108638731
$ws.Range("M25").Value = 1.056713972495097
$ws.Range("N25").Value = 1.016369731102625
